# Apply "New crime data collected" update to the weekly CompStat report
# (24th Precinct, cs-en-us-024pct.xlsx).
#
# 1) Header text: bulletin volume number and the "week covering" date range
#    both advance by one week.
# 2) The data table (rows 15-30, columns C:N) is refreshed with the new
#    week's complaint counts / 28-day counts / percent changes, etc. A few
#    cells that used to be blank placeholders ("0" / "***.*" text) now hold
#    real numbers, so their number format is set explicitly to match the
#    surrounding numeric cells in that column (#,##0 for counts, the
#    parenthesised #,##0.0 format for percentages) before the value is
#    written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header strings -------------------------------------------------

$ws.Range("A8").Value = "Volume 30   Number  44"
$ws.Range("C9").Value = "Report Covering the Week  10/30/2023  Through  11/5/2023"

# ---- Data table (rows 15-30) ----------------------------------------

# --- Row 15 ---
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("D15").Value = 2
$ws.Range("E15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E15").Value = -100
$ws.Range("G15").NumberFormat = "#,##0"
$ws.Range("G15").Value = 2
$ws.Range("H15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H15").Value = -50
$ws.Range("J15").Value = 8
$ws.Range("K15").Value = 25
$ws.Range("N15").Value = -68.75
# --- Row 16 ---
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -80
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 109
$ws.Range("J16").Value = 141
$ws.Range("K16").Value = -22.695035460992
$ws.Range("L16").Value = -9.917355371900
$ws.Range("M16").Value = -35.119047619047
$ws.Range("N16").Value = -85.751633986928
# --- Row 17 ---
$ws.Range("F17").Value = 18
$ws.Range("H17").Value = 50
$ws.Range("I17").Value = 158
$ws.Range("J17").Value = 155
$ws.Range("K17").Value = 1.935483870967
$ws.Range("L17").Value = 3.267973856209
$ws.Range("M17").Value = 83.720930232558
$ws.Range("N17").Value = -62.470308788598
# --- Row 18 ---
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 150
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = 25
$ws.Range("I18").Value = 154
$ws.Range("J18").Value = 156
$ws.Range("K18").Value = -1.282051282051
$ws.Range("L18").Value = 35.087719298245
$ws.Range("M18").Value = 38.738738738738
$ws.Range("N18").Value = -83.242655059847
# --- Row 19 ---
$ws.Range("F19").Value = 60
$ws.Range("G19").Value = 39
$ws.Range("H19").Value = 53.846153846153
$ws.Range("I19").Value = 450
$ws.Range("J19").Value = 448
$ws.Range("K19").Value = 0.446428571428
$ws.Range("L19").Value = 13.065326633165
$ws.Range("M19").Value = -2.173913043478
$ws.Range("N19").Value = -52.330508474576
# --- Row 20 ---
$ws.Range("F20").Value = 5
$ws.Range("H20").Value = -16.666666666666
$ws.Range("I20").Value = 85
$ws.Range("J20").Value = 79
$ws.Range("K20").Value = 7.594936708860
$ws.Range("L20").Value = 11.842105263157
$ws.Range("M20").Value = 142.857142857143
$ws.Range("N20").Value = -89.976415094339
# --- Row 21 ---
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = 13.636363636363
$ws.Range("F21").Value = 110
$ws.Range("G21").Value = 82
$ws.Range("H21").Value = 34.146341463414
$ws.Range("I21").Value = 968
$ws.Range("J21").Value = 988
$ws.Range("K21").Value = -2.024291497975
$ws.Range("L21").Value = 10.628571428571
$ws.Range("M21").Value = 11.009174311926
$ws.Range("N21").Value = -75.493670886075
# --- Row 22 ---
$ws.Range("F22").Value = 5
$ws.Range("G22").Value = 7
$ws.Range("H22").Value = -28.571428571428
$ws.Range("I22").Value = 21
$ws.Range("J22").Value = 27
$ws.Range("K22").Value = -22.222222222222
$ws.Range("L22").Value = -12.5
$ws.Range("M22").Value = -16
# --- Row 23 ---
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 50
$ws.Range("G23").Value = 11
$ws.Range("H23").Value = 9.090909090909
$ws.Range("I23").Value = 99
$ws.Range("J23").Value = 108
$ws.Range("K23").Value = -8.333333333333
$ws.Range("L23").Value = 5.319148936170
$ws.Range("M23").Value = 52.307692307692
# --- Row 24 ---
$ws.Range("C24").Value = 46
$ws.Range("D24").Value = 38
$ws.Range("E24").Value = 21.052631578947
$ws.Range("F24").Value = 120
$ws.Range("G24").Value = 176
$ws.Range("H24").Value = -31.818181818181
$ws.Range("I24").Value = 1345
$ws.Range("J24").Value = 1600
$ws.Range("K24").Value = -15.9375
$ws.Range("L24").Value = 32.642998027613
$ws.Range("M24").Value = 46.994535519125
# --- Row 25 ---
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 125
$ws.Range("F25").Value = 31
$ws.Range("G25").Value = 24
$ws.Range("H25").Value = 29.166666666666
$ws.Range("I25").Value = 259
$ws.Range("J25").Value = 261
$ws.Range("K25").Value = -0.766283524904
$ws.Range("L25").Value = 12.121212121212
$ws.Range("M25").Value = -9.756097560975
# --- Row 26 ---
$ws.Range("D26").NumberFormat = "#,##0"
$ws.Range("D26").Value = 3
$ws.Range("E26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E26").Value = -100
$ws.Range("G26").NumberFormat = "#,##0"
$ws.Range("G26").Value = 3
$ws.Range("H26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H26").Value = -66.666666666666
$ws.Range("J26").Value = 13
$ws.Range("K26").Value = 30.769230769230
# --- Row 27 ---
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("C27").Value = 2
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("D27").Value = 2
$ws.Range("E27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E27").Value = 0
$ws.Range("I27").Value = 38
$ws.Range("J27").Value = 46
$ws.Range("K27").Value = -17.391304347826
$ws.Range("L27").Value = -19.148936170212
# --- Row 28 ---
$ws.Range("L28").Value = -22.222222222222
# --- Row 29 ---
$ws.Range("L29").Value = 20
# --- Row 30 ---
$ws.Range("F30").NumberFormat = "#,##0"
$ws.Range("F30").Value = 3
$ws.Range("I30").Value = 14
$ws.Range("K30").Value = -12.5
$ws.Range("L30").Value = 55.555555555555
